$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("User Details")

for ($r = 11; $r -le 15; $r++) {
    $ws.Cells.Item($r, 1).Value = "user112"
    $ws.Cells.Item($r, 2).Value = "giri"
    $ws.Cells.Item($r, 3).Value = "t"
    $ws.Cells.Item($r, 4).Value = "passwor"
}
